# Rename transcript speaker labels in the DataSheet.
# "HILLARY LEWIS-WOLFSEN" -> "T" (Teacher)
# "STUDENT" -> "S" (Student)
# Only column D (Speaker) is affected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $val = $cell.Value2
    if ($val -eq "HILLARY LEWIS-WOLFSEN") {
        $cell.Value = "T"
    } elseif ($val -eq "STUDENT") {
        $cell.Value = "S"
    }
}
